# V. 111 "EL juego bonito"
# Adds a new movie row ("El juego bonito") to the "Películas" table, just
# above the existing "Alimañas" row, shifting the remaining rows down by
# one and growing the table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")

# Insert a new blank row at row 125 (within the table body); this pushes
# the former rows 125-131 down to 126-132, copying formulas/number formats
# from the surrounding rows.
$ws.Rows.Item(125).Insert()

# Grow the "Tabla24" table/autofilter range to cover the new row.
$tbl = $ws.ListObjects.Item("Tabla24")
$tbl.Resize($ws.Range("B2:I132"))

# Fill in the new movie's data in row 125.
$ws.Range("B125").Value = "El juego bonito"
$ws.Range("C125").Formula = "=AVERAGE(D125,E125,E125,F125,G125,H125,H125,I125)"
$ws.Range("D125").Value = 4
$ws.Range("E125").Value = 3
$ws.Range("F125").Value = 3
$ws.Range("G125").Value = 4
$ws.Range("H125").Value = 6.3
$ws.Range("I125").Value = 4.8

# Re-apply the alignment on the (now shifted-away) "El niño que domó el
# viento" row so it normalizes back to the sheet's common left-aligned
# style used everywhere else in the column.
$ws.Range("B77").HorizontalAlignment = -4131

# Leave the selection where the author ended up, on the new last row.
$ws.Range("C132").Select()
